# "update label form input"
# - "semester list" sheet: column B held the raw semester number (3, 5, 7, ...).
#   That number moves to a new column C, and column B is turned into a
#   Roman-numeral label (III, V, VII, ...) for display. The option-builder
#   formula in column F is updated so the value='' attribute now reads the
#   numeric column C while the visible option text keeps reading column B
#   (which now shows the Roman numeral).
# - Bring "semester list" to the front (it was tab 4/"matkul" before).
# - Refresh the remembered selections on both "semester list" and "matkul".

$wb = $excel.ActiveWorkbook
$wsSemester = $wb.Worksheets.Item("semester list")
$wsMatkul = $wb.Worksheets.Item("matkul")

# Arabic semester number -> Roman numeral label, per row (1-25).
$romans = @{
    1  = "I"
    3  = "III"
    4  = "IV"
    5  = "V"
    6  = "VI"
    7  = "VII"
    9  = "IX"
}

for ($r = 1; $r -le 25; $r++) {
    $num = $wsSemester.Range("B$r").Value2
    $wsSemester.Range("C$r").Value = $num
    $wsSemester.Range("B$r").Value = $romans[[int]$num]
}

# Row 1's F formula isn't part of the shared-formula group (F2:F25 is),
# so update it on its own, then update the rest of the group in one shot
# so Excel keeps them as a shared formula.
$wsSemester.Range("F1").Formula = "=""<option value='""&C1&D1&""' class='""&D1&""'>""&B1&""</option>"""
$wsSemester.Range("F2:F25").Formula = "=""<option value='""&C2&D2&""' class='""&D2&""'>""&B2&""</option>"""

# Bring "semester list" to the foreground and leave a fresh selection on it.
$wsSemester.Activate()
$wsSemester.Range("H17").Select()

# "matkul" keeps its scroll position but gets a new remembered selection
# now that it is no longer the active tab.
$wsMatkul.Range("D137").Select()
$wsSemester.Activate()
